$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 87 (shifts existing rows 87-93 down to 88-94),
# inheriting the formatting of the row above, as Excel does for a normal
# row insert.
$ws.Rows.Item(87).Insert()

# Populate the new row 87 with the "Versicherungs eGK Nummer" variable.
$ws.Range("A87").Value = "Versicherungs eGK Nummer"
$ws.Range("B87").Value = "versicherung_egk"
$ws.Range("C87").Value = "A123456789"
$ws.Range("E87").Value = "regex=[A-Z]\d{9}"

# D/F/H/J/L reuse text that already exists elsewhere on the sheet (the
# "Versicherungsidentifikation" row, now at row 89 after the insert above).
# Copy+paste-values from there instead of re-typing, so plain text like
# "False"/"1..1" doesn't get auto-coerced into a boolean/number.
$ws.Range("D89").Copy()
$ws.Range("D87").PasteSpecial(-4163)
$ws.Range("F89").Copy()
$ws.Range("F87").PasteSpecial(-4163)
$ws.Range("H89").Copy()
$ws.Range("H87").PasteSpecial(-4163)
$ws.Range("J89").Copy()
$ws.Range("J87").PasteSpecial(-4163)
$ws.Range("L89").Copy()
$ws.Range("L87").PasteSpecial(-4163)

# The inserted-row template carried G/I/K/M cells over from the row above;
# the new row doesn't use those columns, so clear them back out.
$ws.Range("K87").ClearContents()
$ws.Range("M87").ClearContents()

# C87 holds free text (an eGK number example) rather than the numeric
# default seen on neighbouring rows, so give it a distinct "text" style:
# the same text number format column C normally uses, but drawn with the
# explicit-black Calibri font used elsewhere in the sheet.
$ws.Range("C87").NumberFormat = "@"
$ws.Range("C87").Font.Color = 0

# Leave the selection/scroll position on the newly added row, as Excel
# would after an in-place insert+edit at that location.
$ws.Range("B85").Select()
